# Webasto ET Plants.xlsx - add missing Plant 106 (Neubrandenburg/Germany)
# row so the filter panel shows the complete plant list.
#
# Strategy: shift the existing data rows (3..22) down by one (working from
# the bottom up so we never overwrite a row before it has been copied),
# then write the new Plant 106 row into row 3 using the same values as
# Plant 101 (row 2) except for the Plant number itself. Formatting is
# reapplied with Copy/PasteSpecial(xlPasteFormats) from an existing data
# row so no brand-new style records get minted - it keeps reusing the
# workbook's existing style indices exactly like the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 22
$lastCol = 7

# 1) Shift rows 3..22 down into rows 4..23, bottom-up so sources are never
#    clobbered before they're read.
for ($r = $lastRow; $r -ge 3; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    for ($c = 1; $c -le $lastCol; $c++) {
        $src = $ws.Cells.Item($srcRow, $c)
        $dst = $ws.Cells.Item($dstRow, $c)
        if ($src.HasFormula) {
            $dst.Formula = $src.Formula
        } else {
            $dst.Value2 = $src.Value2
        }
    }
}

# The newly created row 23 has no formatting yet (it previously didn't
# exist) - copy the formatting from row 2 so it reuses the same style
# records as the rest of the table instead of minting new ones.
$ws.Range("A2:G2").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 2) Write the new Plant 106 row into row 3 - same City/Country/PLZ/
#    Latitude/Longitude as Plant 101, only the Plant number differs.
$ws.Cells.Item(3, 1).Value2 = 106
$ws.Cells.Item(3, 2).Value2 = $ws.Cells.Item(2, 2).Value2
$ws.Cells.Item(3, 3).Value2 = $ws.Cells.Item(2, 3).Value2
$ws.Cells.Item(3, 4).Value2 = $ws.Cells.Item(2, 4).Value2
$ws.Cells.Item(3, 5).Value2 = $ws.Cells.Item(2, 5).Value2
$ws.Cells.Item(3, 6).Value2 = $ws.Cells.Item(2, 6).Value2
$ws.Cells.Item(3, 7).Value2 = $ws.Cells.Item(2, 7).Value2

# Row 3 also needs the standard data-row formatting (it currently still
# carries whatever the old row-3/Plant-131 row had, which is already the
# same style, but make it explicit/robust by reapplying from row 2).
$ws.Range("A2:G2").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Restore the active cell/selection noted in the authored workbook.
$ws.Range("K13").Select() | Out-Null
